# Update activity stats (runs/balls/fours/sixes) on the "till excel form" sheet.
# Values are stored as text (number-stored-as-text) in the original workbook,
# so we prefix with a leading apostrophe to force Excel to keep them as text
# instead of auto-converting to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'12"
$ws.Range("D2").Value = "'11"
$ws.Range("E2").Value = "'1"

$ws.Range("C3").Value = "'58"
$ws.Range("D3").Value = "'29"
$ws.Range("E3").Value = "'8"
$ws.Range("F3").Value = "'2"

$ws.Range("D4").Value = "'2"

$ws.Range("C5").Value = "'0"
$ws.Range("D5").Value = "'3"

$ws.Range("C7").Value = "'6"
$ws.Range("D7").Value = "'8"
$ws.Range("E7").Value = "'0"

$ws.Range("C8").Value = "'1"
$ws.Range("D8").Value = "'3"
$ws.Range("E8").Value = "'0"
$ws.Range("F8").Value = "'0"
